# Remove the leftover "9. " placeholder paragraph (with its trailing
# single-space run) that sits right after the "8. )" code-block line and
# right before the "Contenu de depart" heading.
$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("9. ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $para = $rng.Paragraphs(1)
    # Delete the whole paragraph, including its paragraph mark, so the
    # surrounding paragraphs merge together exactly as in the diff.
    $para.Range.Delete()
}
